$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "Il sistema autentica l’utente`r") {
        # Insert a new list paragraph right after this one, matching its
        # style/numbering (Paragrafoelenco, numId 4), then fill in its text.
        $p.Range.InsertParagraphAfter()
        $p1 = $p.Next()
        $p1.Range.Text = "Il sistema setta le informazioni dell’utente (email, indirizzo e ID)"

        # Insert a second new list paragraph after the first new one.
        $p1.Range.InsertParagraphAfter()
        $p2 = $p1.Next()
        $p2.Range.Text = "Il sistema carica la pagina da mostrare"

        break
    }
}
